$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14: set A14 = 1, change B14 from "green.png" to "orange.png"
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "orange.png"

# Row 20: set A20 = 1, change B20 from "green.png" to "red.png"
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "red.png"

# Row 42: set A42 = 1, change B42 from "green.png" to "red.png"
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "red.png"

# Update the selection on the active sheet to A14
$ws.Range("A14").Select()
